# Apply updated cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so that numeric-looking
# strings (e.g. "230.99", "0.758") are preserved as text, matching
# the original inlineStr cell contents instead of being parsed as numbers.
$cells = @(
  'D2',
  'E2',
  'D3',
  'E3',
  'E4',
  'D5',
  'E5',
  'E6',
  'E7',
  'D8',
  'E8',
  'E9',
  'E10',
  'E11',
  'E12',
  'D13',
  'E13',
  'D14',
  'E14',
  'D15',
  'E15',
  'D16',
  'E16',
  'D17',
  'E17',
  'D18',
  'E18',
  'D19',
  'E19',
  'D20',
  'E20',
  'D21',
  'E21',
  'D22',
  'E22',
  'E23',
  'D24',
  'E24',
  'D25',
  'E25',
  'D26',
  'E26',
  'D27',
  'E27',
  'E28',
  'D29',
  'E29',
  'E30',
  'E31',
  'D32',
  'E32',
  'D33',
  'E33',
  'D34',
  'E34',
  'E35',
  'E36',
  'D37',
  'E37',
  'D38',
  'E38',
  'E39',
  'E40',
  'D41',
  'E41',
  'D42',
  'E42',
  'E43',
  'D44',
  'E44',
  'E45',
  'D46',
  'E46',
  'B47',
  'C47',
  'D47',
  'E47',
  'B48',
  'C48',
  'D48',
  'E48',
  'D49',
  'E49',
  'D50',
  'E50',
  'D51',
  'E51'
)
foreach ($ref in $cells) {
  $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.381.29'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').Value = '2.052.56'
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '230.99'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '57.07'
$ws.Range('E8').Value = '  -3.79%  '
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').Value = '14.64'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('D14').Value = '20.66'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '0.758'
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('D16').Value = '5.29'
$ws.Range('E16').Value = '  -1.41%  '
$ws.Range('D17').Value = '2.056.00'
$ws.Range('E17').Value = '  +0.87%  '
$ws.Range('D18').Value = '37.297.17'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').Value = '6.10'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('D20').Value = '69.66'
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('D21').Value = '0.0₃0824'
$ws.Range('E21').Value = '  -3.07%  '
$ws.Range('D22').Value = '226.30'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '2.40'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').Value = '2.33'
$ws.Range('E25').Value = '  -3.33%  '
$ws.Range('D26').Value = '9.88'
$ws.Range('E26').Value = '  +7.99%  '
$ws.Range('D27').Value = '170.08'
$ws.Range('E27').Value = '  -0.82%  '
$ws.Range('E28').Value = '  -5.90%  '
$ws.Range('D29').Value = '19.21'
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('E30').Value = '  -5.13%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').Value = '4.53'
$ws.Range('E32').Value = '  -4.08%  '
$ws.Range('D33').Value = '0.0622'
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('D34').Value = '4.57'
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('E35').Value = '  -1.49%  '
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').Value = '3.28'
$ws.Range('E37').Value = '  -4.39%  '
$ws.Range('D38').Value = '0.997'
$ws.Range('E38').Value = '  -0.41%  '
$ws.Range('E39').Value = '  -1.89%  '
$ws.Range('E40').Value = '  +2.88%  '
$ws.Range('D41').Value = '98.35'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').Value = '0.0954'
$ws.Range('E42').Value = '  -3.04%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').Value = '1.477.88'
$ws.Range('E44').Value = '  +2.36%  '
$ws.Range('E45').Value = '  +2.53%  '
$ws.Range('D46').Value = '16.63'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '1.03'
$ws.Range('E47').Value = '  -2.92%  '
$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').Value = '3.97'
$ws.Range('E48').Value = '  -5.03%  '
$ws.Range('D49').Value = '7.26'
$ws.Range('E49').Value = '  -1.75%  '
$ws.Range('D50').Value = '2.94'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('D51').Value = '2.239.12'
$ws.Range('E51').Value = '  -1.57%  '
